$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update header row (row 1) column names to short machine-friendly codes
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# 2) Title-case the connector words ("de" -> "De", "del" -> "Del", "y" -> "Y", "el" -> "El")
#    in several municipality / state names for consistency.
$ws.Range("B3").Value = "Pabellón De Arteaga"
$ws.Range("B4").Value = "Rincón De Romos"
$ws.Range("B23").Value = "Guadalupe Y Calvo"
$ws.Range("B25").Value = "Hidalgo Del Parral"
$ws.Range("B40").Value = "Valle De Zaragoza"
$ws.Range("A49").Value = "Ciudad De México"
$ws.Range("A73").Value = "Estado De México"
$ws.Range("B75").Value = "Ecatepec De Morelos"
$ws.Range("B76").Value = "Naucalpan De Juárez"
$ws.Range("B78").Value = "San Felipe Del Progreso"
$ws.Range("B79").Value = "Tenango Del Valle"
$ws.Range("B81").Value = "Tlalnepantla De Baz"
$ws.Range("B84").Value = "Villa De Allende"
$ws.Range("B85").Value = "Villa Del Carbón"
$ws.Range("B91").Value = "Purísima Del Rincón"
$ws.Range("B93").Value = "Silao De La Victoria"
$ws.Range("B96").Value = "Ajuchitlán Del Progreso"
$ws.Range("B97").Value = "Chilapa De Álvarez"
$ws.Range("B99").Value = "Coyuca De Benítez"
$ws.Range("B100").Value = "Coyuca De Catalán"
$ws.Range("B108").Value = "Pachuca De Soto"
$ws.Range("B118").Value = "Lagos De Moreno"
$ws.Range("B119").Value = "Ojuelos De Jalisco"
$ws.Range("B120").Value = "San Juan De Los Lagos"
$ws.Range("B122").Value = "San Miguel El Alto"
$ws.Range("B124").Value = "Tamazula De Gordiano"
$ws.Range("B127").Value = "Unión De Tula"
$ws.Range("B130").Value = "Zapotlán El Grande"
$ws.Range("B150").Value = "Tetela Del Volcán"
$ws.Range("B154").Value = "Ixtlán Del Río"
$ws.Range("B159").Value = "Miahuatlán De Porfirio Díaz"
$ws.Range("B160").Value = "Ocotlán De Morelos"
$ws.Range("B169").Value = "Tlacolula De Matamoros"
$ws.Range("B182").Value = "Tlacotepec De Benito Juárez"
$ws.Range("B189").Value = "Pinal De Amoles"
$ws.Range("B190").Value = "San Juan Del Río"
$ws.Range("B193").Value = "Mexquitic De Carmona"
$ws.Range("B196").Value = "Santa María Del Río"
$ws.Range("B219").Value = "Contla De Juan Cuamatzi"
$ws.Range("B223").Value = "Muñoz De Domingo Arenas"
$ws.Range("B251").Value = "Tlaltenango De Sánchez Román"
$ws.Range("B252").Value = "Trinidad García De La Cadena"
$ws.Range("B253").Value = "Villa De Cos"

# 3) Delete the trailing metadata/footer rows (rows 259-263), which shifts the
#    worksheet dimension from A1:D263 down to A1:D257.
$ws.Range("A259:D263").EntireRow.Delete()

Write-Host "Edit complete."
